$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-style header cells E3:G3 to match A3 (bold, default alignment) and drop H3 entirely
$ws.Range("A3").Copy()
$ws.Range("E3:G3").PasteSpecial(-4122)
$ws.Range("H3").Clear()

# Fill in the new row for the 120x120 benchmark
$ws.Range("A8").Value = "linalg.matmul ins(%arg0, %3 : memref<120x120xf32>, memref<120x16xf32>) outs(%alloc_1 : memref<120x16xf32>)"
$ws.Range("B8").Value = "4.792 " + [char]0x03BC + "s"
$ws.Range("E8").Value = $ws.Range("E9").Value()
$ws.Range("F8").Value = "./bambu-ac_types-clang16.AppImage -v4 --print-dot -lm --soft-float --compiler=I386_CLANG16 --device-name=xcu280-2Lfsvh2892-VVD --clock-period=4 --experimental-setup=BAMBU-BALANCED-MP --channels-number=2 --memory-allocation-policy=ALL_BRAM --disable-function-proxy --generate-tb=forward_kernel_test.xml --simulate --evaluation --simulator=VERILATOR --top-fname=forward_kernel input.ll 2>&1 | tee bambu-log"

# Update F9's Bambu command to the same, newer, single-line command text
$ws.Range("F9").Value = $ws.Range("F8").Value()

$ws.Range("A1").Select()
